# Actualiza la base de datos del Estado de Cuenta:
# se eliminan periodos anteriores para RAFAEL DE JESUS MENDOZA SALCEDO (2412, 2501)
# y se agregan periodos nuevos (2407, 2408); las demas filas se reordenan/actualizan
# para que el periodo de mora quede agrupado y ordenado por trabajador.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$datos = @(
    @{ Row = 16; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2501"; Mora = 52000 },
    @{ Row = 17; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2412"; Mora = 52000 },
    @{ Row = 18; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2411"; Mora = 52000 },
    @{ Row = 19; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2410"; Mora = 52000 },
    @{ Row = 20; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2409"; Mora = 52000 },
    @{ Row = 21; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2408"; Mora = 52000 },
    @{ Row = 22; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2407"; Mora = 52000 },
    @{ Row = 23; Doc = "1052733536"; Nombre = 'LILIANA MARGARITA PRIMO CABEZA'; Periodo = "2501"; Mora = 52000 },
    @{ Row = 24; Doc = "1052733536"; Nombre = 'LILIANA MARGARITA PRIMO CABEZA'; Periodo = "2412"; Mora = 52000 },
    @{ Row = 25; Doc = "1052733536"; Nombre = 'LILIANA MARGARITA PRIMO CABEZA'; Periodo = "2411"; Mora = 52000 },
    @{ Row = 26; Doc = "1052733536"; Nombre = 'LILIANA MARGARITA PRIMO CABEZA'; Periodo = "2410"; Mora = 52000 },
    @{ Row = 27; Doc = "1052733536"; Nombre = 'LILIANA MARGARITA PRIMO CABEZA'; Periodo = "2409"; Mora = 52000 },
    @{ Row = 28; Doc = "1052740154"; Nombre = 'DALI MILENA POLO MENDOZA'; Periodo = "2501"; Mora = 52000 },
    @{ Row = 29; Doc = "1052740154"; Nombre = 'DALI MILENA POLO MENDOZA'; Periodo = "2412"; Mora = 52000 },
    @{ Row = 30; Doc = "1052740154"; Nombre = 'DALI MILENA POLO MENDOZA'; Periodo = "2411"; Mora = 52000 },
    @{ Row = 31; Doc = "1052740154"; Nombre = 'DALI MILENA POLO MENDOZA'; Periodo = "2410"; Mora = 52000 },
    @{ Row = 32; Doc = "1052740154"; Nombre = 'DALI MILENA POLO MENDOZA'; Periodo = "2409"; Mora = 52000 },
    @{ Row = 33; Doc = "1235042119"; Nombre = 'MARIA JOSE VEGA TORDECILLA'; Periodo = "2501"; Mora = 52000 },
    @{ Row = 34; Doc = "1235042119"; Nombre = 'MARIA JOSE VEGA TORDECILLA'; Periodo = "2412"; Mora = 52000 },
    @{ Row = 35; Doc = "1235042119"; Nombre = 'MARIA JOSE VEGA TORDECILLA'; Periodo = "2411"; Mora = 52000 },
    @{ Row = 36; Doc = "1235042119"; Nombre = 'MARIA JOSE VEGA TORDECILLA'; Periodo = "2410"; Mora = 52000 },
    @{ Row = 37; Doc = "1235042119"; Nombre = 'MARIA JOSE VEGA TORDECILLA'; Periodo = "2409"; Mora = 52000 },
    @{ Row = 38; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2503"; Mora = 31200 },
    @{ Row = 39; Doc = "73133162"; Nombre = 'RAFAEL DE JESUS MENDOZA SALCEDO'; Periodo = "2502"; Mora = 52000 },
    @{ Row = 40; Doc = "1052733536"; Nombre = 'LILIANA MARGARITA PRIMO CABEZA'; Periodo = "2503"; Mora = 31200 },
    @{ Row = 41; Doc = "1052733536"; Nombre = 'LILIANA MARGARITA PRIMO CABEZA'; Periodo = "2502"; Mora = 52000 },
    @{ Row = 42; Doc = "1052740154"; Nombre = 'DALI MILENA POLO MENDOZA'; Periodo = "2503"; Mora = 31200 },
    @{ Row = 43; Doc = "1052740154"; Nombre = 'DALI MILENA POLO MENDOZA'; Periodo = "2502"; Mora = 52000 },
    @{ Row = 44; Doc = "1235042119"; Nombre = 'MARIA JOSE VEGA TORDECILLA'; Periodo = "2503"; Mora = 31200 },
    @{ Row = 45; Doc = "1235042119"; Nombre = 'MARIA JOSE VEGA TORDECILLA'; Periodo = "2502"; Mora = 52000 }
)

foreach ($fila in $datos) {
    $r = $fila.Row
    $ws.Range("C$r").Value = $fila.Doc
    $ws.Range("D$r").Value = $fila.Nombre
    $ws.Range("E$r").Value = $fila.Periodo
    $ws.Range("F$r").Value = $fila.Mora
}
